# Update "想去人数" (column F) figures on the 展览 and 全部类型 sheets
# to match the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

$updates = @{
    5  = 495
    6  = 1452
    7  = 771
    9  = 194
    10 = 141
    11 = 193
    12 = 110
    13 = 166
    14 = 149
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
